# Bug fix: printCSV used the wrong sheet name (should read the print
# directory / result-file paths that now live on Sheet3, pointing at the
# new "D:/source/nodejs" location instead of the old "C:/nodejs" one).
$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# Update the two path values on Sheet3.
$ws3.Range("A2").Value = "D:/source/nodejs/tempdata/result.csv"
$ws3.Range("B2").Value = "D:/source/nodejs/print"

# Widen the two columns to fit the longer paths.
$ws3.Columns.Item(1).ColumnWidth = 39.857142857142858
$ws3.Columns.Item(2).ColumnWidth = 24.285714285714285

# Move the active selection on Sheet3 to B6.
$ws3.Activate()
[void]$ws3.Range("B6").Select()
$wb.Worksheets.Item(1).Activate()
